# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型"
# sheets to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 50
$ws1.Range("F4").Value = 3549
$ws1.Range("F5").Value = 2204
$ws1.Range("F9").Value = 61
$ws1.Range("F10").Value = 1306
$ws1.Range("F12").Value = 1815

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 50
$ws4.Range("F4").Value = 3549
$ws4.Range("F5").Value = 2204
$ws4.Range("F10").Value = 61
$ws4.Range("F13").Value = 1306
$ws4.Range("F15").Value = 1815
